$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the topic for cm002 (row 4, column D) with the new text.
$ws.Range("D4").Value = "Grammar of graphics and version control software"

# Leave the final selection on D5, matching the recorded cursor position.
$ws.Range("D5").Select()
